# Add new columns I (I0) and J (IF) to the worksheet, populating header
# and per-row values as described by the commit "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the other header cells (e.g. H1) by
# copying its formatting (bold font + thin border + center/top align)
# onto the two new header cells without disturbing their values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data values for I2:J39 ---
$iValues = @(6,7,8,6,6,7,7,7,6,8,7,8,7,7,8,6,7,6,6,6,5,9,8,8,7,8,7,2,8,8,9,9,8,7,9,5,4,5)
$jValues = @(6,7,8,6,6,7,7,7,6,8,7,8,7,7,8,6,7,6,6,6,5,9,8,8,7,8,7,3,8,8,9,9,8,7,9,5,4,5)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
